$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range('D2') '62.086.14'
$ws.Range('E2').Value = '  -0.20%  '

# Row 3
Set-TextValue $ws.Range('D3') '3.429.24'
$ws.Range('E3').Value = '  -0.03%  '

# Row 4
$ws.Range('E4').Value = '  +0.33%  '

# Row 5
Set-TextValue $ws.Range('D5') '409.19'

# Row 6
Set-TextValue $ws.Range('D6') '130.23'
$ws.Range('E6').Value = '  -1.80%  '

# Row 8
$ws.Range('E8').Value = '  -0.06%  '

# Row 9
$ws.Range('E9').Value = '  +7.10%  '

# Row 10
Set-TextValue $ws.Range('D10') '0.141'
$ws.Range('E10').Value = '  +5.65%  '

# Row 11
Set-TextValue $ws.Range('D11') '42.87'
$ws.Range('E11').Value = '  +1.97%  '

# Row 12
Set-TextValue $ws.Range('D12') '0.0000225'
$ws.Range('E12').Value = '  +53.46%  '

# Row 13
Set-TextValue $ws.Range('D13') '9.27'
$ws.Range('E13').Value = '  +10.53%  '

# Row 14
$ws.Range('E14').Value = '  -0.20%  '

# Row 15
$ws.Range('B15').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C15').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
Set-TextValue $ws.Range('D15') '3.974.16'
$ws.Range('E15').Value = '  +0.08%  '

# Row 16
$ws.Range('B16').Value = 'Chainlink'
$ws.Range('C16').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue $ws.Range('D16') '21.36'
$ws.Range('E16').Value = '  +7.55%  '

# Row 17
Set-TextValue $ws.Range('D17') '3.440.19'
$ws.Range('E17').Value = '  +0.59%  '

# Row 18
Set-TextValue $ws.Range('D18') '12.57'
$ws.Range('E18').Value = '  +7.91%  '

# Row 19
$ws.Range('E19').Value = '  +7.94%  '

# Row 20
Set-TextValue $ws.Range('D20') '62.059.47'
$ws.Range('E20').Value = '  -0.27%  '

# Row 21
Set-TextValue $ws.Range('D21') '456.83'
$ws.Range('E21').Value = '  +46.47%  '

# Row 22
Set-TextValue $ws.Range('D22') '91.69'
$ws.Range('E22').Value = '  +9.27%  '

# Row 23
$ws.Range('E23').Value = '  +1.89%  '

# Row 24
Set-TextValue $ws.Range('D24') '13.11'
$ws.Range('E24').Value = '  +2.61%  '

# Row 25
Set-TextValue $ws.Range('D25') '3.26'
$ws.Range('E25').Value = '  +2.53%  '

# Row 26
Set-TextValue $ws.Range('D26') '33.10'
$ws.Range('E26').Value = '  +11.51%  '

# Row 27
Set-TextValue $ws.Range('D27') '9.08'
$ws.Range('E27').Value = '  +11.28%  '

# Row 28
Set-TextValue $ws.Range('D28') '4.78'
$ws.Range('E28').Value = '  +1.15%  '

# Row 29
$ws.Range('E29').Value = '  -0.91%  '

# Row 30
$ws.Range('B30').Value = 'Cosmos'
$ws.Range('C30').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue $ws.Range('D30') '12.10'
$ws.Range('E30').Value = '  +6.52%  '

# Row 31
$ws.Range('B31').Value = 'Toncoin'
$ws.Range('C31').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue $ws.Range('D31') '2.69'
$ws.Range('E31').Value = '  -2.48%  '

# Row 32
$ws.Range('E32').Value = '  -0.94%  '

# Row 33
$ws.Range('B33').Value = 'InjectiveProtocol'
$ws.Range('C33').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue $ws.Range('D33') '43.05'
$ws.Range('E33').Value = '  -1.59%  '

# Row 34
$ws.Range('B34').Value = 'Hedera'
$ws.Range('C34').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue $ws.Range('D34') '0.114'
$ws.Range('E34').Value = '  -0.66%  '

# Row 35
$ws.Range('E35').Value = '  -0.08%  '

# Row 36
Set-TextValue $ws.Range('D36') '0.0505'
$ws.Range('E36').Value = '  +3.83%  '

# Row 37
Set-TextValue $ws.Range('D37') '54.51'
$ws.Range('E37').Value = '  +5.42%  '

# Row 38
$ws.Range('E38').Value = '  +0.01%  '

# Row 39
Set-TextValue $ws.Range('D39') '3.39'
$ws.Range('E39').Value = '  +1.85%  '

# Row 40
$ws.Range('E40').Value = '  +7.63%  '

# Row 41
$ws.Range('B41').Value = 'Stacks'
$ws.Range('C41').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue $ws.Range('D41') '2.95'
$ws.Range('E41').Value = '  -1.70%  '

# Row 42
$ws.Range('B42').Value = 'TheGraph'
$ws.Range('C42').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
Set-TextValue $ws.Range('D42') '0.320'
$ws.Range('E42').Value = '  +0.62%  '

# Row 43
Set-TextValue $ws.Range('D43') '142.67'
$ws.Range('E43').Value = '  -1.04%  '

# Row 44
Set-TextValue $ws.Range('D44') '4.25'
$ws.Range('E44').Value = '  +9.06%  '

# Row 45
$ws.Range('E45').Value = '  +1.04%  '

# Row 46
$ws.Range('E46').Value = '  +13.36%  '

# Row 47
Set-TextValue $ws.Range('D47') '16.68'
$ws.Range('E47').Value = '  -0.66%  '

# Row 48
Set-TextValue $ws.Range('D48') '22.40'
$ws.Range('E48').Value = '  +5.78%  '

# Row 49
Set-TextValue $ws.Range('D49') '2.14'
$ws.Range('E49').Value = '  +9.57%  '

# Row 50
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue $ws.Range('D50') '0.139'
$ws.Range('E50').Value = '  +16.77%  '

# Row 51
$ws.Range('B51').Value = 'RocketPoolETH'
$ws.Range('C51').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
Set-TextValue $ws.Range('D51') '3.777.44'
$ws.Range('E51').Value = '  -0.24%  '
